$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E to text format before writing, so numeric-looking
# strings (e.g. "0.9967", "985.00") are stored verbatim as text instead of
# being coerced into numbers by the COM Value setter.
$ws.Range("D2:E51").NumberFormat = "@"

# --- Rows whose Coin/Link swapped position (re-ranked) plus new Price/Volume ---
# Row 21 <-> Row 22: WrappedliquidstakedEther2.0 <-> Dai
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "0.9997"
$ws.Range("E21").Value = "  +0.05%  "

$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "2.152.50"
$ws.Range("E22").Value = "  +0.51%  "

# Row 46 <-> Row 47: Quant <-> Aptos
$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").Value = "7.597"
$ws.Range("E46").Value = "  -0.32%  "

$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").Value = "100.95"
$ws.Range("E47").Value = "  +1.09%  "

# --- Price / Volume(1h) refresh for all other rows ---
$ws.Range("D2").Value = "29.962.73"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").Value = "1.906.62"
$ws.Range("E3").Value = "  +0.51%  "
$ws.Range("D4").Value = "0.9967"
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").Value = "0.8384"
$ws.Range("E5").Value = "  +10.27%  "
$ws.Range("D6").Value = "241.78"
$ws.Range("E6").Value = "  +0.83%  "
$ws.Range("D7").Value = "0.9982"
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "0.3196"
$ws.Range("E8").Value = "  +4.47%  "
$ws.Range("D9").Value = "26.73"
$ws.Range("E9").Value = "  +4.67%  "
$ws.Range("D10").Value = "0.06999"
$ws.Range("E10").Value = "  +2.39%  "
$ws.Range("D11").Value = "0.08013"
$ws.Range("E11").Value = "  +0.81%  "
$ws.Range("D12").Value = "0.7521"
$ws.Range("E12").Value = "  +1.15%  "
$ws.Range("D13").Value = "1.916.06"
$ws.Range("E13").Value = "  +1.07%  "
$ws.Range("D14").Value = "5.204"
$ws.Range("E14").Value = "  +0.75%  "
$ws.Range("D15").Value = "92.78"
$ws.Range("E15").Value = "  +2.25%  "
$ws.Range("D16").Value = "29.976.19"
$ws.Range("E16").Value = "  +0.37%  "
$ws.Range("D17").Value = "14.12"
$ws.Range("E17").Value = "  +1.26%  "
$ws.Range("D18").Value = "5.885"
$ws.Range("E18").Value = "  -1.18%  "
$ws.Range("D19").Value = "245.48"
$ws.Range("E19").Value = "  +1.28%  "
$ws.Range("D20").Value = "0.000007760"
$ws.Range("E20").Value = "  +0.96%  "
$ws.Range("D23").Value = "0.9962"
$ws.Range("E23").Value = "  -0.34%  "
$ws.Range("D24").Value = "6.975"
$ws.Range("E24").Value = "  +0.65%  "
$ws.Range("D25").Value = "0.1631"
$ws.Range("E25").Value = "  +26.70%  "
$ws.Range("D26").Value = "169.08"
$ws.Range("E26").Value = "  +1.61%  "
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("D28").Value = "18.90"
$ws.Range("E28").Value = "  +1.13%  "
$ws.Range("D29").Value = "2.090"
$ws.Range("E29").Value = "  +3.01%  "
$ws.Range("E30").Value = "  -1.71%  "
$ws.Range("D31").Value = "1.513"
$ws.Range("E31").Value = "  -0.09%  "
$ws.Range("E32").Value = "  +0.99%  "
$ws.Range("D33").Value = "0.05572"
$ws.Range("E33").Value = "  +6.99%  "
$ws.Range("D34").Value = "4.085"
$ws.Range("E34").Value = "  +0.75%  "
$ws.Range("D35").Value = "1.274"
$ws.Range("E35").Value = "  +1.49%  "
$ws.Range("D36").Value = "0.7342"
$ws.Range("E36").Value = "  +1.17%  "
$ws.Range("D37").Value = "2.706"
$ws.Range("E37").Value = "  -0.15%  "
$ws.Range("D38").Value = "0.01923"
$ws.Range("E38").Value = "  +0.15%  "
$ws.Range("D39").Value = "2.787"
$ws.Range("E39").Value = "  +0.52%  "
$ws.Range("D40").Value = "0.4428"
$ws.Range("E40").Value = "  +0.48%  "
$ws.Range("D41").Value = "72.28"
$ws.Range("E41").Value = "  +0.64%  "
$ws.Range("D42").Value = "5.986"
$ws.Range("E42").Value = "  -2.69%  "
$ws.Range("D43").Value = "0.9981"
$ws.Range("D44").Value = "0.8373"
$ws.Range("E44").Value = "  +1.00%  "
$ws.Range("D45").Value = "1.895"
$ws.Range("E45").Value = "  +0.37%  "
$ws.Range("D48").Value = "9.723"
$ws.Range("E48").Value = "  -0.93%  "
$ws.Range("D49").Value = "985.00"
$ws.Range("E49").Value = "  +9.39%  "
$ws.Range("D50").Value = "2.063.32"
$ws.Range("E50").Value = "  +1.03%  "
$ws.Range("D51").Value = "36.27"
$ws.Range("E51").Value = "  +0.69%  "

# Restore default (no explicit number format) on the touched cells so the
# saved styles match the original sheet (only the header row keeps style s=1).
$ws.Range("D2:E51").ClearFormats()
